$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the BMI Value (C) and BMI Category (D) data for rows 2-8,
# leaving the header row (row 1) intact.
$ws.Range("C2:D8").ClearContents()

# Update the active selection to C2 (was C9).
$ws.Range("C2").Select()
